$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 now carries the label that used to belong to row 6, with new counts.
$ws.Range("A5").Value = "Cilindros Hidráulicos De Alta Pressão"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# Row 6 gets a brand-new label, counts stay the same (1, 1).
$ws.Range("A6").Value = "Cilindros Hidráulicos 700 Bar"
